# Report card table update: refresh indicator rows with 2024 data
# (new commercial revenue year label, rewritten implication text, and
# re-ordered status/time-series values) as described in the commit
# "add edits to report card doc".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New text is introduced in this order so the workbook's shared-string
# table is built up the same way the original author's edits produced it.
$ws.Range("A4").Value = "Commercial revenue (2024 USD)"
$ws.Range("C3").Value = "Number of commercial vessels has been steadily decreasing since around 2000 consistent with decreasing fleet diversity and continued risk to fishery resilience (MAFMC FID)."
$ws.Range("C4").Value = "Average Longfin ex-vessel prices in 2024 increased slightly from 2023 (+10%), but commercial revenue has decreased from 2023 which is most likely driven by a an overall decrease in landings by 23% (MAFMC FID)."
$ws.Range("C2").Value = "An increase in landings since 2020 but decrease in number of vessels could indicate targeted trips in specific times of year and fishers targeting other species when longfin are not available. High variability in landings is common for squid fisheries, and 2024 commercial landings fall within the long term mean."

# Row 2 - Commercial landings
$ws.Range("A2").Value = "Commercial landings (millions of lbs.)`n"
$ws.Range("B2").Value = "Near long term average"
$ws.Range("D2").Value = "Commercial_LONGFINSQUID_Landings_LBS_2025-04-07.png"

# Row 3 - Number of commercial vessels
$ws.Range("A3").Value = "Number of commercial vessels (#)"
$ws.Range("B3").Value = "Below long term average"
$ws.Range("D3").Value = "N_Commercial_Vessels_Landing_LONGFINSQUID_2025-04-07.png"

# Row 4 - Commercial revenue
$ws.Range("B4").Value = "Below long term average"
$ws.Range("D4").Value = "TOTALANNUALREV_LONGFINSQUID_2023Dols_2025-04-07.png"

# Row 5 - Western Gulf Stream Index
$ws.Range("A5").Value = "Western Gulf Stream Index (shift in the western part of the Gulf Stream North wall: mean position: >0 = more northerly, <0 = more southerly)"
$ws.Range("B5").Value = "Above long term average"
$ws.Range("C5").Value = "Since the mid-1990s, north and westward shifts in the Gulf Stream have resulted in an increase in warm core rings and deep water, high salinity heat waves. The position of the Gulf Stream influences seasonal temperature and water mass mixing dynamics that affect longfin squid habitat suitability, temperature-dependent growth, and prey availability."
$ws.Range("D5").Value = "western gulf stream index_2025-04-07.png"

# Row 6 - Bottom temperature
$ws.Range("A6").Value = "Bottom temperature in MAB and SNE(°C) "
$ws.Range("B6").Value = "Above long term average (Fall); near long term average (Spring)"
$ws.Range("C6").Value = "Longfin squid seasonal distribution and growth rates are likely temperature dependent, avoiding water <8°C. Inshore temperature thresholds (around 14°C) initiate migration of squid from offshore overwintering habitats. "
$ws.Range("D6").Value = "BottomT_2025-04-17.png"

# Match the author's final selection in the saved worksheet view
$ws.Range("C2").Select()
